# Apply corrections to "Mortendad" -> "Mortandad" and restructure the
# "Los Alamos and Pajarito Canyons" section into two separate sections
# on the "Regional Exhibit" sheet.

$wb = $excel.ActiveWorkbook

$wsMap = $wb.Worksheets.Item("Regional for Mapping")
$wsExh = $wb.Worksheets.Item("Regional Exhibit")

# --- Sheet "Regional for Mapping": fix misspelling in Watershed column (P) ---
$wsMap.Range("P3").Value = "Mortandad"
$wsMap.Range("P4").Value = "Mortandad"
$wsMap.Range("P5").Value = "Mortandad"
$wsMap.Range("P6").Value = "Mortandad"

# --- Sheet "Regional Exhibit" ---

# Fix misspelling in the section header for Mortandad Canyon
$wsExh.Range("A5").Value = "Mortandad Canyon"

# Rename combined section header to refer only to Los Alamos Canyon
$wsExh.Range("A10").Value = "Los Alamos Canyon"

# Insert a new section header row for Pajarito Canyon before the wells
# that were previously lumped under "Los Alamos and Pajarito Canyons"
$wsExh.Rows.Item(17).Insert()
$wsExh.Range("A17").Value = "Pajarito Canyon"
$wsExh.Range("A17:H17").Merge()

# Copy the formatting of an existing section header row onto the new row
$wsExh.Range("A10:H10").Copy()
$wsExh.Range("A17:H17").PasteSpecial(-4122)

# Widen column C slightly to accommodate the data
$wsExh.Columns.Item(3).ColumnWidth = 12.14

$excel.CutCopyMode = $false
